$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: mark the "ned" run as a detected grammar error too, matching
#         the existing spellStart/spellEnd proofErr pair already there
#         (adds <w:proofErr w:type="gramStart"/> / "gramEnd" around it).
#         We rebuild paragraph 2 verbatim (same runs/rsids) via InsertXML
#         so only the two proofErr markers are actually new.
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs(2)
$para2Xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="002E1C50" w:rsidRDefault="002E1C50">' +
            '<w:r><w:t xml:space="preserve">Student of </w:t></w:r>' +
            '<w:proofErr w:type="spellStart"/>' +
            '<w:proofErr w:type="gramStart"/>' +
            '<w:r><w:t>ned</w:t></w:r>' +
            '<w:proofErr w:type="spellEnd"/>' +
            '<w:proofErr w:type="gramEnd"/>' +
            '<w:r><w:t xml:space="preserve"> university 1</w:t></w:r>' +
            '<w:r w:rsidRPr="002E1C50"><w:rPr><w:vertAlign w:val="superscript"/></w:rPr><w:t>st</w:t></w:r>' +
            '<w:r><w:t xml:space="preserve"> year </w:t></w:r>' +
            '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
            '</w:p>'
$p2.Range.InsertXML($para2Xml)

# ---------------------------------------------------------------------
# Step 2: split the paragraph right after "year " into its own new
#         paragraph (the bookmark, sitting at that boundary, travels
#         into the new paragraph automatically).
# ---------------------------------------------------------------------
$d.Content.Find.Execute("year ", $true, $false, $false, $false, $false, $true, 1, $false, "year ^p", 2)

# ---------------------------------------------------------------------
# Step 3: add the submission-date text, split across the bookmark so
#         the final layout is: "Submission date 24" <bookmark> "-10-2024"
#         Insert the text that follows the bookmark first (keeps the
#         bookmark pinned at its current slot), then insert the text
#         that precedes the bookmark (pushes the bookmark after it).
# ---------------------------------------------------------------------
$bm = $d.Bookmarks("_GoBack")
$afterBookmark = $d.Range($bm.End, $bm.End)
$afterBookmark.InsertAfter("-10-2024")

$bm2 = $d.Bookmarks("_GoBack")
$beforeBookmark = $d.Range($bm2.Start, $bm2.Start)
$beforeBookmark.InsertBefore("Submission date 24")
